$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting the existing rows 47-49 down to 48-50.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly record.
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44714
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = 100112037
$ws.Range("G47").Value = "Cebollín"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 180
$ws.Range("K47").Value = 6500
$ws.Range("L47").Value = 7000
$ws.Range("M47").Value = 6778
$ws.Range("N47").Value = "$/paquete 36 unidades"
$ws.Range("O47").Value = "Región Metropolitana"
$ws.Range("P47").Value = 188
$ws.Range("Q47").Value = 36
$ws.Range("R47").Value = "Hortaliza"
